$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Standard MS failure Notes")

# New column D header + 21 data rows (rows 1-22), referencing new shared strings 155-176
$values = @(
    "R replacement",
    "MS Okay note",
    "MS note okay note",
    "Wobble okay note",
    "Wobble failure",
    "NA",
    "Stellaris Failures",
    "Baseswap",
    "No full length product",
    "Wrong Mass found",
    "Poor Fluoresence ",
    "Plate scrambled",
    "Poor Purity",
    "Impurity at specific mass",
    "General Synthesis Failure",
    "Poor Texas Red coupling",
    "Poor Pulsar Coupling",
    "Poor Methlyene Blue Coupling",
    "poor Quasar 705 Coupling",
    "Poor HEX coupling",
    "Poor TET coupling",
    "Poor TAM coupling"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

# Column width changes: C narrower, new D column wider
$ws.Columns.Item(3).ColumnWidth = 20.665
$ws.Columns.Item(4).ColumnWidth = 61.5

# Move the active selection to C11
$ws.Range("C11").Select()
